{"js": "// The \"on_delete\" explanation paragraph is rewritten:\n//   - a new run is inserted right after the \"on_delete   \" label run,\n//   - \"\u0443\u043a\u0430\u0437\u044b\u0432\u0430\u0435\u0442\u044c\" -> \"\u0441\u0441\u044b\u043b\u0430\u0435\u0442\u044c\u0441\u044f\" (still spell-flagged, so a\n//     <w:proofErr w:type=\"spellStart\"/>...<w:proofErr w:type=\"spellEnd\"/>\n//     pair keeps wrapping it),\n//   - the old tail \" \u043a\u0430\u043a \u0431\u0443\u0434\u0443\u0442 \u0443\u0434\u0430\u043b\u044f\u0442\u044c\u0441\u044f \u0437\u0430\u043f\u0438\u0441\u0438  \u0432 \u043e\u0431\u043e\u0445 \u0442\u0430\u0431\u043b\u0438\u0446\u0430\u0445\" (plus its\n//     proofErr wrappers for the \"\u0437\u0430\u043f\u0438\u0441\u0438  \u0432\" grammar flag and the \"\u043e\u0431\u043e\u0445\"\n//     spelling flag) is replaced by \" \u0442\u0435\u043a\u0443\u0449\u0430\u044f\".\n//\n// A plain text search+delete strips run text but leaves the stand-alone\n// <w:proofErr/> markers behind (they live *between* runs, not inside one),\n// so the reliable way to reproduce the exact target markup is to rebuild\n// the paragraph's trailing OOXML in one shot via insertOoxml(...,\n// Word.InsertLocation.replace). The untouched leading part of the\n// paragraph (its w14:paraId/rsid attributes, <w:pPr>, and the bold\n// \"on_delete   \" label runs) is read back from getOoxml() so nothing about\n// it has to be hard-coded.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"on_delete\") === 0 && text.indexOf(\"\u0443\u0434\u0430\u043b\u044f\u0442\u044c\u0441\u044f\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the on_delete paragraph to edit.\");\n}\n\nconst ooxmlResult = target.getOoxml();\nawait context.sync();\nconst xml = ooxmlResult.value;\n\nconst bodyIdx = xml.indexOf(\"<w:body>\");\nconst pStart = xml.indexOf(\"<w:p\", bodyIdx);\nconst pEnd = xml.indexOf(\"</w:p>\", pStart) + \"</w:p>\".length;\nconst fullParagraphXml = xml.substring(pStart, pEnd);\n\n// Keep everything up to (and including) the run that holds the three\n// spaces right after the bold \"on_delete\" label untouched.\nconst labelMarker = '<w:t xml:space=\"preserve\">   </w:t></w:r>';\nconst labelEnd = fullParagraphXml.indexOf(labelMarker);\nif (labelEnd === -1) {\n  throw new Error(\"Could not locate the on_delete label run boundary.\");\n}\nconst headerXml = fullParagraphXml.substring(\n  0,\n  labelEnd + labelMarker.length\n);\n\nconst newTailXml =\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\">\u0443\u043a\u0430\u0437\u044b\u0432\u0430\u0435\u0442 \u0447\u0442\u043e \u0431\u0443\u0434\u0435\u0442 \u043d\u0430 \u043c\u0435\u0441\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u0432 \u0442\u0435\u043a\u0443\u0449\u0435\u0439 \u0442\u0430\u0431\u043b\u0438\u0446\u0435, \u0435\u0441\u043b\u0438 \u0443\u0434\u0430\u043b\u044f\u0435\u043c \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u0432\u043e \u0432\u0442\u043e\u0440\u043e\u0439 \u0442\u0430\u0431\u043b\u0438\u0446\u0435, \u043d\u0430 \u043a\u043e\u0442\u043e\u0440\u0443\u044e </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>\u0441\u0441\u044b\u043b\u0430\u0435\u0442\u044c\u0441\u044f</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> \u0442\u0435\u043a\u0443\u0449\u0430\u044f</w:t></w:r>';\n\nconst newParagraphXml = headerXml + newTailXml + \"</w:p>\";\n\nconst packageXml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  newParagraphXml +\n  \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\ntarget.getRange().insertOoxml(packageXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"on_delete\" explanation paragraph is rewritten:\n#   - a new run is inserted right after the \"on_delete   \" label run,\n#   - \"\u0443\u043a\u0430\u0437\u044b\u0432\u0430\u0435\u0442\u044c\" -> \"\u0441\u0441\u044b\u043b\u0430\u0435\u0442\u044c\u0441\u044f\" (still spell-flagged, so a\n#     <w:proofErr w:type=\"spellStart\"/>...<w:proofErr w:type=\"spellEnd\"/>\n#     pair keeps wrapping it),\n#   - the old tail \" \u043a\u0430\u043a \u0431\u0443\u0434\u0443\u0442 \u0443\u0434\u0430\u043b\u044f\u0442\u044c\u0441\u044f \u0437\u0430\u043f\u0438\u0441\u0438  \u0432 \u043e\u0431\u043e\u0445 \u0442\u0430\u0431\u043b\u0438\u0446\u0430\u0445\" (plus its\n#     proofErr wrappers for the \"\u0437\u0430\u043f\u0438\u0441\u0438  \u0432\" grammar flag and the \"\u043e\u0431\u043e\u0445\"\n#     spelling flag) is replaced by \" \u0442\u0435\u043a\u0443\u0449\u0430\u044f\".\n#\n# A plain Find/Replace (or Range.Text assignment) on the descriptive tail\n# strips the run text but leaves the stand-alone <w:proofErr/> markers\n# behind (they live *between* runs, not inside one), so the reliable way to\n# reproduce the exact target markup is to rebuild the paragraph's trailing\n# OOXML in one shot via Range.InsertXML(...). The untouched leading part of\n# the paragraph (its w14:paraId/rsid attributes, <w:pPr>, and the bold\n# \"on_delete   \" label runs) is read back from Range.WordOpenXML so nothing\n# about it has to be hard-coded.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($target -eq $null) {\n        $t = $p.Range.Text\n        if ($t.IndexOf(\"on_delete\") -eq 0 -and $t.IndexOf(\"\u0443\u0434\u0430\u043b\u044f\u0442\u044c\u0441\u044f\") -ge 0) {\n            $target = $p\n        }\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the on_delete paragraph to edit.\"\n}\n\n$r = $target.Range\n$xml = $r.WordOpenXML\n\n$bodyIdx = $xml.IndexOf(\"<w:body>\")\n$pStart = $xml.IndexOf(\"<w:p\", $bodyIdx)\n$pEndTag = \"</w:p>\"\n$pEnd = $xml.IndexOf($pEndTag, $pStart) + $pEndTag.Length\n$fullParagraphXml = $xml.Substring($pStart, $pEnd - $pStart)\n\n# Keep everything up to (and including) the run that holds the three spaces\n# right after the bold \"on_delete\" label untouched.\n$labelMarker = '<w:t xml:space=\"preserve\">   </w:t></w:r>'\n$labelEnd = $fullParagraphXml.IndexOf($labelMarker)\nif ($labelEnd -eq -1) {\n    throw \"Could not locate the on_delete label run boundary.\"\n}\n$headerXml = $fullParagraphXml.Substring(0, $labelEnd + $labelMarker.Length)\n\n$newTailXml = '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n    '<w:t xml:space=\"preserve\">\u0443\u043a\u0430\u0437\u044b\u0432\u0430\u0435\u0442 \u0447\u0442\u043e \u0431\u0443\u0434\u0435\u0442 \u043d\u0430 \u043c\u0435\u0441\u0442\u043e \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u0432 \u0442\u0435\u043a\u0443\u0449\u0435\u0439 \u0442\u0430\u0431\u043b\u0438\u0446\u0435, \u0435\u0441\u043b\u0438 \u0443\u0434\u0430\u043b\u044f\u0435\u043c \u0437\u043d\u0430\u0447\u0435\u043d\u0438\u0435 \u0432\u043e \u0432\u0442\u043e\u0440\u043e\u0439 \u0442\u0430\u0431\u043b\u0438\u0446\u0435, \u043d\u0430 \u043a\u043e\u0442\u043e\u0440\u0443\u044e </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t>\u0441\u0441\u044b\u043b\u0430\u0435\u0442\u044c\u0441\u044f</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> \u0442\u0435\u043a\u0443\u0449\u0430\u044f</w:t></w:r>'\n\n$newParagraphXml = $headerXml + $newTailXml + \"</w:p>\"\n\n$packageXml = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    $newParagraphXml +\n    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$r.InsertXML($packageXml)\n"}
